# Remove column from alcohol measurement data:
# Column M (a duplicate/obsolete measurement column) is deleted from
# Sheet1. Deleting the column shifts the old column N one position to
# the left so it becomes the new column M, and the used range shrinks
# from A1:N119 to A1:M119.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns("M").Delete() | Out-Null

# Leave the active cell where the deleted column used to be, matching
# the post-edit selection state.
$ws.Range("M1").Select() | Out-Null
